$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "15+14=29"
$t.Cell(1, 2).Range.Text = "6+46=52"
$t.Cell(1, 3).Range.Text = "78+4=82"
$t.Cell(1, 4).Range.Text = "84-78=6"
$t.Cell(1, 5).Range.Text = "8+44=52"
$t.Cell(2, 1).Range.Text = "88-63=25"
$t.Cell(2, 2).Range.Text = "6+13=19"
$t.Cell(2, 3).Range.Text = "29+61=90"
$t.Cell(2, 4).Range.Text = "80+8=88"
$t.Cell(2, 5).Range.Text = "8+26=34"
$t.Cell(3, 1).Range.Text = "95-27=68"
$t.Cell(3, 2).Range.Text = "70-34=36"
$t.Cell(3, 3).Range.Text = "70-49=21"
$t.Cell(3, 4).Range.Text = "43+10=53"
$t.Cell(3, 5).Range.Text = "72-17=55"
$t.Cell(4, 1).Range.Text = "68-7=61"
$t.Cell(4, 2).Range.Text = "14+2=16"
$t.Cell(4, 3).Range.Text = "22+68=90"
$t.Cell(4, 4).Range.Text = "49-1=48"
$t.Cell(4, 5).Range.Text = "71-7=64"
$t.Cell(5, 1).Range.Text = "27-8=19"
$t.Cell(5, 2).Range.Text = "51+46=97"
$t.Cell(5, 3).Range.Text = "82-46=36"
$t.Cell(5, 4).Range.Text = "87-82=5"
$t.Cell(5, 5).Range.Text = "94-85=9"
$t.Cell(6, 1).Range.Text = "34+13=47"
$t.Cell(6, 2).Range.Text = "58+0=58"
$t.Cell(6, 3).Range.Text = "23+24=47"
$t.Cell(6, 4).Range.Text = "24+43=67"
$t.Cell(6, 5).Range.Text = "99-46=53"
$t.Cell(7, 1).Range.Text = "23-1=22"
$t.Cell(7, 2).Range.Text = "37+56=93"
$t.Cell(7, 3).Range.Text = "4+42=46"
$t.Cell(7, 4).Range.Text = "53+5=58"
$t.Cell(7, 5).Range.Text = "33-31=2"
$t.Cell(8, 1).Range.Text = "47-25=22"
$t.Cell(8, 2).Range.Text = "26+58=84"
$t.Cell(8, 3).Range.Text = "14+54=68"
$t.Cell(8, 4).Range.Text = "28+11=39"
$t.Cell(8, 5).Range.Text = "64-15=49"
$t.Cell(9, 1).Range.Text = "0+35=35"
$t.Cell(9, 2).Range.Text = "86+2=88"
$t.Cell(9, 3).Range.Text = "52-34=18"
$t.Cell(9, 4).Range.Text = "20+59=79"
$t.Cell(9, 5).Range.Text = "11+30=41"
$t.Cell(10, 1).Range.Text = "13+57=70"
$t.Cell(10, 2).Range.Text = "13+12=25"
$t.Cell(10, 3).Range.Text = "57+38=95"
$t.Cell(10, 4).Range.Text = "88-76=12"
$t.Cell(10, 5).Range.Text = "20+4=24"
$t.Cell(11, 1).Range.Text = "43+13=56"
$t.Cell(11, 2).Range.Text = "20+73=93"
$t.Cell(11, 3).Range.Text = "41+48=89"
$t.Cell(11, 4).Range.Text = "17-9=8"
$t.Cell(11, 5).Range.Text = "73-27=46"
$t.Cell(12, 1).Range.Text = "22+71=93"
$t.Cell(12, 2).Range.Text = "79+13=92"
$t.Cell(12, 3).Range.Text = "11+33=44"
$t.Cell(12, 4).Range.Text = "68-36=32"
$t.Cell(12, 5).Range.Text = "80-49=31"
$t.Cell(13, 1).Range.Text = "26+41=67"
$t.Cell(13, 2).Range.Text = "71-23=48"
$t.Cell(13, 3).Range.Text = "8+68=76"
$t.Cell(13, 4).Range.Text = "0+19=19"
$t.Cell(13, 5).Range.Text = "34+20=54"
$t.Cell(14, 1).Range.Text = "93-85=8"
$t.Cell(14, 2).Range.Text = "19+51=70"
$t.Cell(14, 3).Range.Text = "18+80=98"
$t.Cell(14, 4).Range.Text = "37+28=65"
$t.Cell(14, 5).Range.Text = "16+83=99"
$t.Cell(15, 1).Range.Text = "1+32=33"
$t.Cell(15, 2).Range.Text = "98-42=56"
$t.Cell(15, 3).Range.Text = "50+44=94"
$t.Cell(15, 4).Range.Text = "56+36=92"
$t.Cell(15, 5).Range.Text = "53+35=88"
$t.Cell(16, 1).Range.Text = "16+14=30"
$t.Cell(16, 2).Range.Text = "57-11=46"
$t.Cell(16, 3).Range.Text = "84-25=59"
$t.Cell(16, 4).Range.Text = "18-5=13"
$t.Cell(16, 5).Range.Text = "31+50=81"
$t.Cell(17, 1).Range.Text = "75+15=90"
$t.Cell(17, 2).Range.Text = "38-10=28"
$t.Cell(17, 3).Range.Text = "90-58=32"
$t.Cell(17, 4).Range.Text = "89-65=24"
$t.Cell(17, 5).Range.Text = "53+33=86"
$t.Cell(18, 1).Range.Text = "90-81=9"
$t.Cell(18, 2).Range.Text = "20+54=74"
$t.Cell(18, 3).Range.Text = "96-66=30"
$t.Cell(18, 4).Range.Text = "87-42=45"
$t.Cell(18, 5).Range.Text = "40-34=6"
$t.Cell(19, 1).Range.Text = "8+58=66"
$t.Cell(19, 2).Range.Text = "18+65=83"
$t.Cell(19, 3).Range.Text = "9+90=99"
$t.Cell(19, 4).Range.Text = "44-5=39"
$t.Cell(19, 5).Range.Text = "98-26=72"
$t.Cell(20, 1).Range.Text = "17+45=62"
$t.Cell(20, 2).Range.Text = "57-21=36"
$t.Cell(20, 3).Range.Text = "1+75=76"
$t.Cell(20, 4).Range.Text = "28+6=34"
$t.Cell(20, 5).Range.Text = "84+9=93"
